$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.022.18"
$ws.Range("E2").Value = "  +0.29%  "

# Row 3
$ws.Range("D3").Value = "1.643.15"
$ws.Range("E3").Value = "  +0.30%  "

# Row 4
$ws.Range("E4").Value = "  +0.28%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.32"
$ws.Range("E5").Value = "  +0.24%  "

# Row 6
$ws.Range("E6").Value = "  -0.02%  "

# Row 7
$ws.Range("E7").Value = "  +0.26%  "

# Row 8
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.255"
$ws.Range("E8").Value = "  +0.15%  "

# Row 9
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0638"
$ws.Range("E9").Value = "  +0.27%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.59"
$ws.Range("E10").Value = "  -0.36%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0796"
$ws.Range("E11").Value = "  +0.19%  "

# Row 12
$ws.Range("E12").Value = "  +0.14%  "

# Row 13
$ws.Range("D13").Value = "1.647.46"
$ws.Range("E13").Value = "  +0.90%  "

# Row 14
$ws.Range("E14").Value = "  -0.13%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.44"
$ws.Range("E15").Value = "  +1.42%  "

# Row 16
$ws.Range("E16").Value = "  +0.46%  "

# Row 17
$ws.Range("D17").Value = "26.051.96"
$ws.Range("E17").Value = "  +0.33%  "

# Row 18
$ws.Range("E18").Value = "  +0.28%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "194.68"

# Row 20
$ws.Range("E20").Value = "  -0.48%  "

# Row 21
$ws.Range("E21").Value = "  -0.35%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.20"
$ws.Range("E22").Value = "  -1.12%  "

# Row 23
$ws.Range("E23").Value = "  +4.39%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "143.99"
$ws.Range("E24").Value = "  -0.18%  "

# Row 25
$ws.Range("E25").Value = "  -0.26%  "

# Row 26
$ws.Range("E26").Value = "  +0.05%  "

# Row 27
$ws.Range("E27").Value = "  +0.52%  "

# Row 28
$ws.Range("E28").Value = "  +0.06%  "

# Row 29
$ws.Range("E29").Value = "  +0.32%  "

# Row 30
$ws.Range("E30").Value = "  -1.05%  "

# Row 31
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.28"
$ws.Range("E31").Value = "  -0.76%  "

# Row 32
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.26"
$ws.Range("E32").Value = "  +0.86%  "

# Row 33
$ws.Range("E33").Value = "  -0.23%  "

# Row 35
$ws.Range("E35").Value = "  +0.02%  "

# Row 36
$ws.Range("D36").Value = "1.130.56"
$ws.Range("E36").Value = "  -0.73%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.538"
$ws.Range("E37").Value = "  -1.55%  "

# Row 38
$ws.Range("E38").Value = "  -0.14%  "

# Row 39
$ws.Range("E39").Value = "  -0.29%  "

# Row 40
$ws.Range("E40").Value = "  +0.43%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.76"
$ws.Range("E41").Value = "  -0.54%  "

# Row 42
$ws.Range("E42").Value = "  -0.63%  "

# Row 43
$ws.Range("E43").Value = "  +0.56%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "56.48"
$ws.Range("E44").Value = "  +0.04%  "

# Row 45
$ws.Range("E45").Value = "  +2.57%  "

# Row 46
$ws.Range("E46").Value = "  -1.56%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.78"
$ws.Range("E47").Value = "  +2.08%  "

# Row 48
$ws.Range("E48").Value = "  -0.25%  "

# Row 49
$ws.Range("E49").Value = "  +0.12%  "

# Row 50
$ws.Range("E50").Value = "  -1.35%  "

# Row 51
$ws.Range("E51").Value = "  -0.16%  "

